{"js": "// Ordered list of [oldText, newText] pairs, one per paragraph in the document\n// (title paragraph first, then each table cell's paragraph in row-major order).\nconst REPLACEMENTS = [\n  [\"2025-02-04 Tuesday\", \"2025-02-05 Wednesday\"],\n  [\"77+17=\", \"95-46=\"],\n  [\"18+4=\", \"8+34=\"],\n  [\"62-49=\", \"42-29=\"],\n  [\"16+27=\", \"93-88=\"],\n  [\"71-9=\", \"95-48=\"],\n  [\"8+8=\", \"71-15=\"],\n  [\"42-15=\", \"76+18=\"],\n  [\"49+22=\", \"56+28=\"],\n  [\"17+5=\", \"94-25=\"],\n  [\"82-38=\", \"35-6=\"],\n  [\"58+33=\", \"8+17=\"],\n  [\"64-26=\", \"29+24=\"],\n  [\"52+9=\", \"71-29=\"],\n  [\"49+49=\", \"56-27=\"],\n  [\"58+37=\", \"29+19=\"],\n  [\"76-48=\", \"82-13=\"],\n  [\"49+19=\", \"45+49=\"],\n  [\"85-78=\", \"14+19=\"],\n  [\"71-47=\", \"83-36=\"],\n  [\"37+56=\", \"13+49=\"],\n  [\"75+9=\", \"27-18=\"],\n  [\"35+47=\", \"48+37=\"],\n  [\"6+75=\", \"92-44=\"],\n  [\"49+14=\", \"82-18=\"],\n  [\"94-26=\", \"55-39=\"],\n  [\"92-53=\", \"25+27=\"],\n  [\"96-39=\", \"67+29=\"],\n  [\"8+56=\", \"40-8=\"],\n  [\"40-6=\", \"92-48=\"],\n  [\"73-39=\", \"37-9=\"],\n  [\"84-76=\", \"38+8=\"],\n  [\"47+8=\", \"19+2=\"],\n  [\"22+9=\", \"87-49=\"],\n  [\"80-46=\", \"96-59=\"],\n  [\"47+15=\", \"90-21=\"],\n  [\"68-59=\", \"83-79=\"],\n  [\"7+49=\", \"72-56=\"],\n  [\"70-44=\", \"46+9=\"],\n  [\"17+35=\", \"42-35=\"],\n  [\"68+4=\", \"51-5=\"],\n  [\"19+39=\", \"17+18=\"],\n  [\"28+49=\", \"32-15=\"],\n  [\"9+42=\", \"60-39=\"],\n  [\"74-37=\", \"88-79=\"],\n  [\"72-27=\", \"88-29=\"],\n  [\"15+46=\", \"81-8=\"],\n  [\"83+8=\", \"72-69=\"],\n  [\"22+49=\", \"38+9=\"],\n  [\"30-19=\", \"53-18=\"],\n  [\"73-26=\", \"4+79=\"],\n  [\"41-33=\", \"49+39=\"],\n  [\"8+63=\", \"55+18=\"],\n  [\"6+65=\", \"51-28=\"],\n  [\"39+48=\", \"59+13=\"],\n  [\"93-38=\", \"19+28=\"],\n  [\"39+37=\", \"7+27=\"],\n  [\"92-43=\", \"71-69=\"],\n  [\"78-69=\", \"60-38=\"],\n  [\"40-15=\", \"49+45=\"],\n  [\"7+66=\", \"72-26=\"],\n  [\"80-24=\", \"23-15=\"],\n  [\"84-26=\", \"19+76=\"],\n  [\"83-55=\", \"44+19=\"],\n  [\"36-7=\", \"14+39=\"],\n  [\"54+28=\", \"15+39=\"],\n  [\"26+48=\", \"14+17=\"],\n  [\"23-19=\", \"24+69=\"],\n  [\"80-36=\", \"55+8=\"],\n  [\"9+77=\", \"37+7=\"],\n  [\"71-26=\", \"48+4=\"],\n  [\"67-18=\", \"47+18=\"],\n  [\"16+25=\", \"26+67=\"],\n  [\"65+7=\", \"14+49=\"],\n  [\"64-5=\", \"93-85=\"],\n  [\"22-13=\", \"91-68=\"],\n  [\"90-41=\", \"68+14=\"],\n  [\"60-34=\", \"35-27=\"],\n  [\"39+5=\", \"67+7=\"],\n  [\"8+39=\", \"39+34=\"],\n  [\"57+36=\", \"7+65=\"],\n  [\"80-62=\", \"31-22=\"],\n  [\"26+9=\", \"3+88=\"],\n  [\"71-34=\", \"7+34=\"],\n  [\"51-29=\", \"9+84=\"],\n  [\"84+8=\", \"8+68=\"],\n  [\"94-46=\", \"6+55=\"],\n  [\"65-7=\", \"26+27=\"],\n  [\"97-29=\", \"18+49=\"],\n  [\"36+29=\", \"77+4=\"],\n  [\"82-4=\", \"8+86=\"],\n  [\"17+25=\", \"8+53=\"],\n  [\"17+28=\", \"80-42=\"],\n  [\"13+79=\", \"61-12=\"],\n  [\"28-9=\", \"51-7=\"],\n  [\"30-21=\", \"80-35=\"],\n  [\"50-24=\", \"74+7=\"],\n  [\"44+38=\", \"66+19=\"],\n  [\"42-8=\", \"9+76=\"],\n  [\"71-13=\", \"94-5=\"],\n  [\"55-9=\", \"14+39=\"],\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nif (items.length !== REPLACEMENTS.length) {\n  throw new Error(\n    \"Paragraph count mismatch: expected \" + REPLACEMENTS.length + \" got \" + items.length\n  );\n}\n\nfor (let i = 0; i < items.length; i++) {\n  const [oldText, newText] = REPLACEMENTS[i];\n  const para = items[i];\n  if (para.text !== oldText) {\n    throw new Error(\n      \"Paragraph \" + i + \" text mismatch: expected \" + JSON.stringify(oldText) +\n      \" got \" + JSON.stringify(para.text)\n    );\n  }\n  if (oldText !== newText) {\n    para.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Ordered list of [oldText, newText] pairs: the title/date line, followed by\n# every arithmetic-expression cell in the table (row-major order). Each old\n# value is unique in the document, so a whole-document Find/Replace for each\n# pair is unambiguous.\n$pairs = @(\n    @(\"2025-02-04 Tuesday\", \"2025-02-05 Wednesday\"),\n    @(\"77+17=\", \"95-46=\"),\n    @(\"18+4=\", \"8+34=\"),\n    @(\"62-49=\", \"42-29=\"),\n    @(\"16+27=\", \"93-88=\"),\n    @(\"71-9=\", \"95-48=\"),\n    @(\"8+8=\", \"71-15=\"),\n    @(\"42-15=\", \"76+18=\"),\n    @(\"49+22=\", \"56+28=\"),\n    @(\"17+5=\", \"94-25=\"),\n    @(\"82-38=\", \"35-6=\"),\n    @(\"58+33=\", \"8+17=\"),\n    @(\"64-26=\", \"29+24=\"),\n    @(\"52+9=\", \"71-29=\"),\n    @(\"49+49=\", \"56-27=\"),\n    @(\"58+37=\", \"29+19=\"),\n    @(\"76-48=\", \"82-13=\"),\n    @(\"49+19=\", \"45+49=\"),\n    @(\"85-78=\", \"14+19=\"),\n    @(\"71-47=\", \"83-36=\"),\n    @(\"37+56=\", \"13+49=\"),\n    @(\"75+9=\", \"27-18=\"),\n    @(\"35+47=\", \"48+37=\"),\n    @(\"6+75=\", \"92-44=\"),\n    @(\"49+14=\", \"82-18=\"),\n    @(\"94-26=\", \"55-39=\"),\n    @(\"92-53=\", \"25+27=\"),\n    @(\"96-39=\", \"67+29=\"),\n    @(\"8+56=\", \"40-8=\"),\n    @(\"40-6=\", \"92-48=\"),\n    @(\"73-39=\", \"37-9=\"),\n    @(\"84-76=\", \"38+8=\"),\n    @(\"47+8=\", \"19+2=\"),\n    @(\"22+9=\", \"87-49=\"),\n    @(\"80-46=\", \"96-59=\"),\n    @(\"47+15=\", \"90-21=\"),\n    @(\"68-59=\", \"83-79=\"),\n    @(\"7+49=\", \"72-56=\"),\n    @(\"70-44=\", \"46+9=\"),\n    @(\"17+35=\", \"42-35=\"),\n    @(\"68+4=\", \"51-5=\"),\n    @(\"19+39=\", \"17+18=\"),\n    @(\"28+49=\", \"32-15=\"),\n    @(\"9+42=\", \"60-39=\"),\n    @(\"74-37=\", \"88-79=\"),\n    @(\"72-27=\", \"88-29=\"),\n    @(\"15+46=\", \"81-8=\"),\n    @(\"83+8=\", \"72-69=\"),\n    @(\"22+49=\", \"38+9=\"),\n    @(\"30-19=\", \"53-18=\"),\n    @(\"73-26=\", \"4+79=\"),\n    @(\"41-33=\", \"49+39=\"),\n    @(\"8+63=\", \"55+18=\"),\n    @(\"6+65=\", \"51-28=\"),\n    @(\"39+48=\", \"59+13=\"),\n    @(\"93-38=\", \"19+28=\"),\n    @(\"39+37=\", \"7+27=\"),\n    @(\"92-43=\", \"71-69=\"),\n    @(\"78-69=\", \"60-38=\"),\n    @(\"40-15=\", \"49+45=\"),\n    @(\"7+66=\", \"72-26=\"),\n    @(\"80-24=\", \"23-15=\"),\n    @(\"84-26=\", \"19+76=\"),\n    @(\"83-55=\", \"44+19=\"),\n    @(\"36-7=\", \"14+39=\"),\n    @(\"54+28=\", \"15+39=\"),\n    @(\"26+48=\", \"14+17=\"),\n    @(\"23-19=\", \"24+69=\"),\n    @(\"80-36=\", \"55+8=\"),\n    @(\"9+77=\", \"37+7=\"),\n    @(\"71-26=\", \"48+4=\"),\n    @(\"67-18=\", \"47+18=\"),\n    @(\"16+25=\", \"26+67=\"),\n    @(\"65+7=\", \"14+49=\"),\n    @(\"64-5=\", \"93-85=\"),\n    @(\"22-13=\", \"91-68=\"),\n    @(\"90-41=\", \"68+14=\"),\n    @(\"60-34=\", \"35-27=\"),\n    @(\"39+5=\", \"67+7=\"),\n    @(\"8+39=\", \"39+34=\"),\n    @(\"57+36=\", \"7+65=\"),\n    @(\"80-62=\", \"31-22=\"),\n    @(\"26+9=\", \"3+88=\"),\n    @(\"71-34=\", \"7+34=\"),\n    @(\"51-29=\", \"9+84=\"),\n    @(\"84+8=\", \"8+68=\"),\n    @(\"94-46=\", \"6+55=\"),\n    @(\"65-7=\", \"26+27=\"),\n    @(\"97-29=\", \"18+49=\"),\n    @(\"36+29=\", \"77+4=\"),\n    @(\"82-4=\", \"8+86=\"),\n    @(\"17+25=\", \"8+53=\"),\n    @(\"17+28=\", \"80-42=\"),\n    @(\"13+79=\", \"61-12=\"),\n    @(\"28-9=\", \"51-7=\"),\n    @(\"30-21=\", \"80-35=\"),\n    @(\"50-24=\", \"74+7=\"),\n    @(\"44+38=\", \"66+19=\"),\n    @(\"42-8=\", \"9+76=\"),\n    @(\"71-13=\", \"94-5=\"),\n    @(\"55-9=\", \"14+39=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    if ($oldText -eq $newText) { continue }\n\n    $range = $d.Content\n    # Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n    #   MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format,\n    #   ReplaceWith, Replace)  -- wdFindContinue=1, wdReplaceAll=2\n    $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Could not find text: $oldText\"\n    }\n}\n\n$d.Save()"}
